$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("test_existing_survey_import 1").Name = "existing_survey_import_1_test"
$wb.Worksheets.Item("test_existing_survey_import 2").Name = "existing_survey_import_2_test"

# Make the second sheet the active sheet (so activeTab=1 and tabSelected moves to sheet 2)
$wb.Worksheets.Item("existing_survey_import_2_test").Activate()
